$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (Dr. Hadas Kress-Gazit talk): fill in real title + abstract (replacing "Title coming soon!" / " " placeholders)
$ws.Range("G31").Value = 'Will Robots Kill Us?'
$ws.Range("H31").Value = 'You know the sci fi trope — humanity depends on robots that decide to kill us all until a hero shows up to stop their evil plan. How close is that to reality? Will robots take over the world? Are they dangerous? Will they be? In this talk I''ll describe what robots can and cannot do these days and discuss what the future might look like.'

# Row 30 (Dr. Maddie Reynolds talk): fill in real title + abstract
$ws.Range("H30").Value = ' New York State prisons are low-information environments, where there is no internet access and books and media are subject to strict censorship. How can librarians on the outside respond to incarcerated patrons'' need for books and information? The Cornell Prison Education Program (CPEP) is grappling with these questions as they seek to create research opportunities and expand information access in prison in anticipation of their launch of a Cornell BA inside. In this talk, we''ll discuss what it''s like to build a library program for incarcerated students, introducing Cornell''s world class research materials behind the walls. We''ll learn about CPEP and its history, nationwide efforts to bring library resources to incarcerated patrons, and what we''re doing at Cornell to provide incarcerated students with more academic resources. We''ll also talk about how to get involved in CPEP! '
$ws.Range("G30").Value = 'Bringing the Library to Prison with the Cornell Prison Education Program'

# H31's abstract is long -> wrap text and expand the row to Excel's max height, like the source edit
$ws.Range("H31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 409.5

# Update the saved view/selection state to reflect scrolling down to the new rows
$ws.Range("D34").Select()

"done"
